$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '72.065.55'
$ws.Range("E2").Value = '  +0.18%  '

$ws.Range("D3").Value = '4.016.48'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '529.40'
$ws.Range("E5").Value = '  +1.18%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '151.26'
$ws.Range("E6").Value = '  +1.38%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.708'
$ws.Range("E7").Value = '  +12.80%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.749'
$ws.Range("E9").Value = '  +1.28%  '

$ws.Range("E10").Value = '  -3.06%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000326'
$ws.Range("E11").Value = '  -4.38%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '48.34'
$ws.Range("E12").Value = '  +4.11%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.62'
$ws.Range("E13").Value = '  -1.25%  '

$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '4.649.27'
$ws.Range("E14").Value = '  -1.11%  '

$ws.Range("D15").Value = '4.018.00'
$ws.Range("E15").Value = '  -1.10%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.03'
$ws.Range("E16").Value = '  -2.18%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '20.47'
$ws.Range("E17").Value = '  -4.86%  '

$ws.Range("E18").Value = '  -1.01%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.19'
$ws.Range("E19").Value = '  -2.58%  '

$ws.Range("D20").Value = '72.069.19'
$ws.Range("E20").Value = '  +0.06%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '428.07'
$ws.Range("E21").Value = '  -3.26%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '97.78'
$ws.Range("E22").Value = '  +2.22%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.48'
$ws.Range("E23").Value = '  -1.72%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.21'
$ws.Range("E24").Value = '  +3.07%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '14.19'
$ws.Range("E25").Value = '  -1.53%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.34'
$ws.Range("E26").Value = '  -8.83%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.68'
$ws.Range("E27").Value = '  -5.25%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.85'
$ws.Range("E28").Value = '  +1.25%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '36.72'
$ws.Range("E29").Value = '  -1.66%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '3.60'
$ws.Range("E30").Value = '  +16.72%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '13.33'
$ws.Range("E31").Value = '  -1.77%  '

$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.130'
$ws.Range("E32").Value = '  -0.08%  '

$ws.Range("B33").Value = 'Bittensor'
$ws.Range("C33").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '681.43'
$ws.Range("E33").Value = '  -3.03%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.15'
$ws.Range("E34").Value = '  +2.93%  '

$ws.Range("B35").Value = 'InjectiveProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '44.80'
$ws.Range("E35").Value = '  +9.08%  '

$ws.Range("B36").Value = 'OKB'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '65.61'
$ws.Range("E36").Value = '  -2.75%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.450'
$ws.Range("E37").Value = '  +0.36%  '

$ws.Range("D38").Value = '0.0₃0827'
$ws.Range("E38").Value = '  -9.56%  '

$ws.Range("B39").Value = 'ThetaToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.44'
$ws.Range("E39").Value = '  -5.37%  '

$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.150'
$ws.Range("E40").Value = '  -2.98%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.07%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.997'
$ws.Range("E42").Value = '  -0.16%  '

$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0486'
$ws.Range("E43").Value = '  -0.99%  '

$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.24'
$ws.Range("E44").Value = '  +3.61%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.150'
$ws.Range("E45").Value = '  +2.34%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.71'
$ws.Range("E46").Value = '  -3.22%  '

$ws.Range("B47").Value = 'THORChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.85'
$ws.Range("E47").Value = '  +7.31%  '

$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.37'
$ws.Range("E48").Value = '  -5.06%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.02'
$ws.Range("E49").Value = '  -5.64%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000270'
$ws.Range("E50").Value = '  -4.53%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '142.32'
$ws.Range("E51").Value = '  -0.87%  '
